$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 116
$ws.Range("I33").Value = 61
$ws.Range("K33").Value = 61
$ws.Range("M33").Value = 168
$ws.Range("H127").Value = 55556348
$ws.Range("I127").Value = 543.7778
$ws.Range("J127").Value = 111112150
$ws.Range("K127").Value = 1631.3334
$ws.Range("L127").Value = 333336450
$ws.Range("M127").Value = 3328.6666
$ws.Range("N127").Value = -333346370
$ws.Range("H132").Value = 3048.8235
$ws.Range("I132").Value = 2614.375
$ws.Range("K132").Value = 7843.125
$ws.Range("M132").Value = -5313.125
$ws.Range("H133").Value = 75000
$ws.Range("J133").Value = 75000
$ws.Range("L133").Value = 75000
$ws.Range("N133").Value = -85120

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1617.5264
$ws.Range("I2").Value = 1617.5264
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1617.5264
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -1504.5264
$ws.Range("H61").Value = 7726.551
$ws.Range("I61").Value = 4514.8125
$ws.Range("J61").Value = 13772.177
$ws.Range("K61").Value = 4514.8125
$ws.Range("L61").Value = 13772.177
$ws.Range("M61").Value = -4302.8125
$ws.Range("N61").Value = -14196.177
$ws.Range("H74").Value = 4593.5264
$ws.Range("I74").Value = 2149.862
$ws.Range("J74").Value = 12467.556
$ws.Range("K74").Value = 2149.862
$ws.Range("L74").Value = 12467.556
$ws.Range("M74").Value = -1275.862
$ws.Range("N74").Value = -14215.556
$ws.Range("H77").Value = 4593.5264
$ws.Range("I77").Value = 2149.862
$ws.Range("J77").Value = 12467.556
$ws.Range("K77").Value = 10749.31
$ws.Range("L77").Value = 62337.78
$ws.Range("M77").Value = -6381.310000000001
$ws.Range("N77").Value = -71073.78
$ws.Range("H88").Value = 9352.286
$ws.Range("I88").Value = 11493.2
$ws.Range("J88").Value = 4000
$ws.Range("K88").Value = 11493.2
$ws.Range("L88").Value = 4000
$ws.Range("M88").Value = -11087.2
$ws.Range("N88").Value = -4812
$ws.Range("H91").Value = 9352.286
$ws.Range("I91").Value = 11493.2
$ws.Range("J91").Value = 4000
$ws.Range("K91").Value = 11493.2
$ws.Range("L91").Value = 4000
$ws.Range("M91").Value = -10089.2
$ws.Range("N91").Value = -6808
$ws.Range("H116").Value = 1617.5264
$ws.Range("I116").Value = 1617.5264
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1617.5264
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = 676.4736
$ws.Range("H122").Value = 8931138
$ws.Range("I122").Value = 2919.6365
$ws.Range("J122").Value = 41667936
$ws.Range("K122").Value = 8758.9095
$ws.Range("L122").Value = 125003808
$ws.Range("M122").Value = -6308.9095
$ws.Range("N122").Value = -125008708
$ws.Range("H136").Value = 7726.551
$ws.Range("I136").Value = 4514.8125
$ws.Range("J136").Value = 13772.177
$ws.Range("K136").Value = 13544.4375
$ws.Range("L136").Value = 41316.531
$ws.Range("M136").Value = -10994.4375
$ws.Range("N136").Value = -46416.531

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1617.5264
$ws.Range("I3").Value = 1617.5264
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1617.5264
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -1503.5264
$ws.Range("H86").Value = 7094065.5
$ws.Range("I86").Value = 9010872
$ws.Range("J86").Value = 1880.7
$ws.Range("K86").Value = 9010872
$ws.Range("L86").Value = 1880.7
$ws.Range("M86").Value = -9009749
$ws.Range("N86").Value = -4126.7
$ws.Range("H89").Value = 7094065.5
$ws.Range("I89").Value = 9010872
$ws.Range("J89").Value = 1880.7
$ws.Range("K89").Value = 45054360
$ws.Range("L89").Value = 9403.5
$ws.Range("M89").Value = -45048744
$ws.Range("N89").Value = -20635.5
$ws.Range("H134").Value = 68247.664
$ws.Range("I134").Value = 1327.3636
$ws.Range("J134").Value = 252278.5
$ws.Range("K134").Value = 3982.0908
$ws.Range("L134").Value = 756835.5
$ws.Range("M134").Value = -1447.0908
$ws.Range("N134").Value = -761905.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 100000
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 100000
$ws.Range("K10").Value = 0
$ws.Range("L10").ClearContents()
$ws.Range("M10").Value = 100000
$ws.Range("N10").Value = -100278
$ws.Range("H31").Value = 1736.878
$ws.Range("I31").Value = 1146.3243
$ws.Range("J31").Value = 7199.5
$ws.Range("K31").Value = 1146.3243
$ws.Range("L31").Value = 7199.5
$ws.Range("M31").Value = -851.3243
$ws.Range("N31").Value = -7789.5
$ws.Range("H34").Value = 1736.878
$ws.Range("I34").Value = 1146.3243
$ws.Range("J34").Value = 7199.5
$ws.Range("K34").Value = 1146.3243
$ws.Range("L34").Value = 7199.5
$ws.Range("M34").Value = -944.3243
$ws.Range("N34").Value = -7603.5
$ws.Range("H58").Value = 2934628.2
$ws.Range("I58").Value = 5052220.5
$ws.Range("J58").Value = 2577.077
$ws.Range("K58").Value = 5052220.5
$ws.Range("L58").Value = 2577.077
$ws.Range("M58").Value = -5052017.5
$ws.Range("N58").Value = -2983.077
$ws.Range("H86").Value = 1516.875
$ws.Range("I86").Value = 1462.8462
$ws.Range("J86").Value = 1751
$ws.Range("K86").Value = 1462.8462
$ws.Range("L86").Value = 1751
$ws.Range("M86").Value = -339.8462
$ws.Range("N86").Value = -3997
$ws.Range("H89").Value = 1516.875
$ws.Range("I89").Value = 1462.8462
$ws.Range("J89").Value = 1751
$ws.Range("K89").Value = 7314.231
$ws.Range("L89").Value = 8755
$ws.Range("M89").Value = -1698.231
$ws.Range("N89").Value = -19987
$ws.Range("H94").Value = 892
$ws.Range("I94").Value = 756
$ws.Range("J94").Value = 908
$ws.Range("K94").Value = 756
$ws.Range("L94").Value = 908
$ws.Range("M94").Value = -305
$ws.Range("N94").Value = -1810
$ws.Range("H99").Value = 1462.9166
$ws.Range("I99").Value = 1038.125
$ws.Range("J99").Value = 2312.5
$ws.Range("K99").Value = 1038.125
$ws.Range("L99").Value = 2312.5
$ws.Range("M99").Value = 459.875
$ws.Range("N99").Value = -5308.5
$ws.Range("H122").Value = 9844.571
$ws.Range("I122").Value = 9844.571
$ws.Range("K122").Value = 29533.713
$ws.Range("M122").Value = -27083.713
$ws.Range("H126").Value = 1462.9166
$ws.Range("I126").Value = 1038.125
$ws.Range("J126").Value = 2312.5
$ws.Range("K126").Value = 3114.375
$ws.Range("L126").Value = 6937.5
$ws.Range("M126").Value = -644.375
$ws.Range("N126").Value = -11877.5
$ws.Range("H132").Value = 2483.0278
$ws.Range("I132").Value = 2429.4443
$ws.Range("J132").Value = 2643.7778
$ws.Range("K132").Value = 7288.3329
$ws.Range("L132").Value = 7931.3334
$ws.Range("M132").Value = -4758.3329
$ws.Range("N132").Value = -12991.3334
$ws.Range("H134").Value = 2613.257
$ws.Range("I134").Value = 2283.3845
$ws.Range("K134").Value = 6850.1535
$ws.Range("M134").Value = -4315.1535
$ws.Range("H136").Value = 2934628.2
$ws.Range("I136").Value = 5052220.5
$ws.Range("J136").Value = 2577.077
$ws.Range("K136").Value = 15156661.5
$ws.Range("L136").Value = 7731.231000000001
$ws.Range("M136").Value = -15154111.5
$ws.Range("N136").Value = -12831.231

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 20833770
$ws.Range("I5").Value = 497.2857
$ws.Range("K5").Value = 1491.8571
$ws.Range("M5").Value = -1379.8571
$ws.Range("H12").Value = 34482970
$ws.Range("I12").Value = 76923250
$ws.Range("J12").Value = 239.5625
$ws.Range("K12").Value = 230769750
$ws.Range("L12").Value = 718.6875
$ws.Range("M12").Value = -230769577
$ws.Range("N12").Value = -1064.6875
$ws.Range("H38").Value = 60.3
$ws.Range("I38").Value = 34.166668
$ws.Range("J38").Value = 99.5
$ws.Range("K38").Value = 102.500004
$ws.Range("L38").Value = 298.5
$ws.Range("M38").Value = 244.499996
$ws.Range("N38").Value = -992.5
$ws.Range("H135").Value = 20833770
$ws.Range("I135").Value = 497.2857
$ws.Range("K135").Value = 4475.571300000001
$ws.Range("M135").Value = -1940.571300000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2417.1738
$ws.Range("I126").Value = 1866.3334
$ws.Range("K126").Value = 5599.0002
$ws.Range("M126").Value = -3129.0002
$ws.Range("H132").Value = 9838.733
$ws.Range("I132").Value = 3488.75
$ws.Range("J132").Value = 12147.818
$ws.Range("K132").Value = 10466.25
$ws.Range("L132").Value = 36443.454
$ws.Range("M132").Value = -7936.25
$ws.Range("N132").Value = -41503.454

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6111.5557
$ws.Range("I7").Value = 5625.5
$ws.Range("J7").Value = 10000
$ws.Range("K7").Value = 5625.5
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = -5513.5
$ws.Range("N7").Value = -10224
$ws.Range("H40").Value = 3586.2144
$ws.Range("I40").Value = 3267.25
$ws.Range("K40").Value = 3267.25
$ws.Range("M40").Value = -3131.25
$ws.Range("H122").Value = 7548.8286
$ws.Range("I122").Value = 7139.5
$ws.Range("J122").Value = 8241.538
$ws.Range("K122").Value = 21418.5
$ws.Range("L122").Value = 24724.614
$ws.Range("M122").Value = -18968.5
$ws.Range("N122").Value = -29624.614
$ws.Range("H126").Value = 6111.5557
$ws.Range("I126").Value = 5625.5
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 16876.5
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = -14406.5
$ws.Range("N126").Value = -34940

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").ClearContents()
$ws.Range("N7").Value = 0
$ws.Range("H113").Value = 777.1739
$ws.Range("I113").Value = 453
$ws.Range("J113").Value = 1074.3334
$ws.Range("K113").Value = 1359
$ws.Range("L113").Value = 3223.0002
$ws.Range("M113").Value = 811
$ws.Range("N113").Value = -7563.0002
$ws.Range("H117").Value = 52366.668
$ws.Range("J117").Value = 52366.668
$ws.Range("L117").Value = 52366.668
$ws.Range("N117").Value = -61544.668
$ws.Range("H126").Value = 1335.0869
$ws.Range("I126").Value = 1319.381
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 3958.143
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -1488.143
$ws.Range("N126").Value = -9440
$ws.Range("H136").Value = 6219.0293
$ws.Range("I136").Value = 2096.9473
$ws.Range("J136").Value = 11440.333
$ws.Range("K136").Value = 6290.841899999999
$ws.Range("L136").Value = 34320.999
$ws.Range("M136").Value = -3740.841899999999
$ws.Range("N136").Value = -39420.999
